$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to match the latest scrape.
# NumberFormat is forced to Text ("@") before assignment so that numeric-looking
# strings (e.g. "1.0000", "316.24") are stored as literal text instead of being
# coerced into numbers by Excel, matching the original inline-string cell layout.
# Style is then reset to "Normal" so no stray style index is left on the cell.

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '24.529.77'
$r.Style = 'Normal'

$r = $ws.Range('E2')
$r.NumberFormat = '@'
$r.Value = '  -0.77%  '
$r.Style = 'Normal'

$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '1.697.23'
$r.Style = 'Normal'

$r = $ws.Range('E3')
$r.NumberFormat = '@'
$r.Value = '  -0.16%  '
$r.Style = 'Normal'

$r = $ws.Range('D4')
$r.NumberFormat = '@'
$r.Value = '1.0000'
$r.Style = 'Normal'

$r = $ws.Range('E4')
$r.NumberFormat = '@'
$r.Value = '  -0.50%  '
$r.Style = 'Normal'

$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '316.24'
$r.Style = 'Normal'

$r = $ws.Range('E5')
$r.NumberFormat = '@'
$r.Value = '  +0.07%  '
$r.Style = 'Normal'

$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '1.000'
$r.Style = 'Normal'

$r = $ws.Range('E6')
$r.NumberFormat = '@'
$r.Value = '  -0.58%  '
$r.Style = 'Normal'

$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '0.3911'
$r.Style = 'Normal'

$r = $ws.Range('E7')
$r.NumberFormat = '@'
$r.Value = '  -0.61%  '
$r.Style = 'Normal'

$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.4079'
$r.Style = 'Normal'

$r = $ws.Range('E8')
$r.NumberFormat = '@'
$r.Value = '  +0.91%  '
$r.Style = 'Normal'

$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '1.495'
$r.Style = 'Normal'

$r = $ws.Range('E9')
$r.NumberFormat = '@'
$r.Value = '  -1.94%  '
$r.Style = 'Normal'

$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.9996'
$r.Style = 'Normal'

$r = $ws.Range('E10')
$r.NumberFormat = '@'
$r.Value = '  -0.56%  '
$r.Style = 'Normal'

$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '52.28'
$r.Style = 'Normal'

$r = $ws.Range('E11')
$r.NumberFormat = '@'
$r.Value = '  -2.69%  '
$r.Style = 'Normal'

$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '0.08823'
$r.Style = 'Normal'

$r = $ws.Range('E12')
$r.NumberFormat = '@'
$r.Value = '  -0.53%  '
$r.Style = 'Normal'

$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '26.68'
$r.Style = 'Normal'

$r = $ws.Range('E13')
$r.NumberFormat = '@'
$r.Value = '  +12.89%  '
$r.Style = 'Normal'

$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '7.536'
$r.Style = 'Normal'

$r = $ws.Range('E14')
$r.NumberFormat = '@'
$r.Value = '  +1.82%  '
$r.Style = 'Normal'

$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '8.189'
$r.Style = 'Normal'

$r = $ws.Range('E15')
$r.NumberFormat = '@'
$r.Value = '  +0.04%  '
$r.Style = 'Normal'

$r = $ws.Range('E16')
$r.NumberFormat = '@'
$r.Value = '  +2.04%  '
$r.Style = 'Normal'

$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '1.688.51'
$r.Style = 'Normal'

$r = $ws.Range('E17')
$r.NumberFormat = '@'
$r.Value = '  -1.30%  '
$r.Style = 'Normal'

$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '98.11'
$r.Style = 'Normal'

$r = $ws.Range('E18')
$r.NumberFormat = '@'
$r.Value = '  -1.61%  '
$r.Style = 'Normal'

$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '0.07176'
$r.Style = 'Normal'

$r = $ws.Range('E19')
$r.NumberFormat = '@'
$r.Value = '  +1.73%  '
$r.Style = 'Normal'

$r = $ws.Range('E20')
$r.NumberFormat = '@'
$r.Value = '  +5.15%  '
$r.Style = 'Normal'

$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '7.320'
$r.Style = 'Normal'

$r = $ws.Range('E21')
$r.NumberFormat = '@'
$r.Value = '  +3.30%  '
$r.Style = 'Normal'

$r = $ws.Range('E22')
$r.NumberFormat = '@'
$r.Value = '  -0.58%  '
$r.Style = 'Normal'

$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '14.39'
$r.Style = 'Normal'

$r = $ws.Range('E23')
$r.NumberFormat = '@'
$r.Value = '  -2.52%  '
$r.Style = 'Normal'

$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '24.524.85'
$r.Style = 'Normal'

$r = $ws.Range('E24')
$r.NumberFormat = '@'
$r.Value = '  -0.79%  '
$r.Style = 'Normal'

$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '3.047'
$r.Style = 'Normal'

$r = $ws.Range('E25')
$r.NumberFormat = '@'
$r.Value = '  -2.70%  '
$r.Style = 'Normal'

$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '2.333'
$r.Style = 'Normal'

$r = $ws.Range('E26')
$r.NumberFormat = '@'
$r.Value = '  -1.52%  '
$r.Style = 'Normal'

$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '22.99'
$r.Style = 'Normal'

$r = $ws.Range('E27')
$r.NumberFormat = '@'
$r.Value = '  +1.04%  '
$r.Style = 'Normal'

$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '167.98'
$r.Style = 'Normal'

$r = $ws.Range('E28')
$r.NumberFormat = '@'
$r.Value = '  +2.95%  '
$r.Style = 'Normal'

$r = $ws.Range('B29')
$r.NumberFormat = '@'
$r.Value = 'BitcoinCash'
$r.Style = 'Normal'

$r = $ws.Range('C29')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$r.Style = 'Normal'

$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '146.78'
$r.Style = 'Normal'

$r = $ws.Range('E29')
$r.NumberFormat = '@'
$r.Value = '  +8.06%  '
$r.Style = 'Normal'

$r = $ws.Range('B30')
$r.NumberFormat = '@'
$r.Value = 'Filecoin'
$r.Style = 'Normal'

$r = $ws.Range('C30')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$r.Style = 'Normal'

$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '8.503'
$r.Style = 'Normal'

$r = $ws.Range('E30')
$r.NumberFormat = '@'
$r.Value = '  -3.00%  '
$r.Style = 'Normal'

$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '5.394'
$r.Style = 'Normal'

$r = $ws.Range('E31')
$r.NumberFormat = '@'
$r.Value = '  +4.12%  '
$r.Style = 'Normal'

$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '2.215'
$r.Style = 'Normal'

$r = $ws.Range('E32')
$r.NumberFormat = '@'
$r.Value = '  +11.80%  '
$r.Style = 'Normal'

$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '1.875.25'
$r.Style = 'Normal'

$r = $ws.Range('E33')
$r.NumberFormat = '@'
$r.Value = '  -1.31%  '
$r.Style = 'Normal'

$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '0.08781'
$r.Style = 'Normal'

$r = $ws.Range('E34')
$r.NumberFormat = '@'
$r.Value = '  -2.87%  '
$r.Style = 'Normal'

$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '7.330'
$r.Style = 'Normal'

$r = $ws.Range('E35')
$r.NumberFormat = '@'
$r.Value = '  -4.32%  '
$r.Style = 'Normal'

$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '1.047'
$r.Style = 'Normal'

$r = $ws.Range('E36')
$r.NumberFormat = '@'
$r.Value = '  -1.90%  '
$r.Style = 'Normal'

$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '0.03059'
$r.Style = 'Normal'

$r = $ws.Range('E37')
$r.NumberFormat = '@'
$r.Value = '  +10.14%  '
$r.Style = 'Normal'

$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '0.2806'
$r.Style = 'Normal'

$r = $ws.Range('E38')
$r.NumberFormat = '@'
$r.Value = '  +1.82%  '
$r.Style = 'Normal'

$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '10.97'
$r.Style = 'Normal'

$r = $ws.Range('E39')
$r.NumberFormat = '@'
$r.Value = '  -1.51%  '
$r.Style = 'Normal'

$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '0.09193'
$r.Style = 'Normal'

$r = $ws.Range('E40')
$r.NumberFormat = '@'
$r.Value = '  +0.41%  '
$r.Style = 'Normal'

$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '14.26'
$r.Style = 'Normal'

$r = $ws.Range('E41')
$r.NumberFormat = '@'
$r.Value = '  -1.81%  '
$r.Style = 'Normal'

$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '0.8044'
$r.Style = 'Normal'

$r = $ws.Range('E42')
$r.NumberFormat = '@'
$r.Value = '  +4.91%  '
$r.Style = 'Normal'

$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '1.482'
$r.Style = 'Normal'

$r = $ws.Range('E43')
$r.NumberFormat = '@'
$r.Value = '  +1.29%  '
$r.Style = 'Normal'

$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '17.46'
$r.Style = 'Normal'

$r = $ws.Range('E44')
$r.NumberFormat = '@'
$r.Value = '  +10.03%  '
$r.Style = 'Normal'

$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '2.681'
$r.Style = 'Normal'

$r = $ws.Range('E45')
$r.NumberFormat = '@'
$r.Value = '  +4.23%  '
$r.Style = 'Normal'

$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '0.7292'
$r.Style = 'Normal'

$r = $ws.Range('E46')
$r.NumberFormat = '@'
$r.Value = '  +1.58%  '
$r.Style = 'Normal'

$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '4.272'
$r.Style = 'Normal'

$r = $ws.Range('E47')
$r.NumberFormat = '@'
$r.Value = '  +1.34%  '
$r.Style = 'Normal'

$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '1.419'
$r.Style = 'Normal'

$r = $ws.Range('E48')
$r.NumberFormat = '@'
$r.Value = '  +6.13%  '
$r.Style = 'Normal'

$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '0.9991'
$r.Style = 'Normal'

$r = $ws.Range('E49')
$r.NumberFormat = '@'
$r.Value = '  -0.53%  '
$r.Style = 'Normal'

$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '141.39'
$r.Style = 'Normal'

$r = $ws.Range('E50')
$r.NumberFormat = '@'
$r.Value = '  +0.88%  '
$r.Style = 'Normal'

$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '0.08161'
$r.Style = 'Normal'

$r = $ws.Range('E51')
$r.NumberFormat = '@'
$r.Value = '  +2.17%  '
$r.Style = 'Normal'

